$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 328210
$ws.Cells.Item(2, 4).Value = 418052235
$ws.Cells.Item(4, 3).Value = 330
$ws.Cells.Item(4, 4).Value = 472192
$ws.Cells.Item(10, 3).Value = 119020
$ws.Cells.Item(10, 4).Value = 174383802
$ws.Cells.Item(12, 3).Value = 60902
$ws.Cells.Item(12, 4).Value = 87886707
$ws.Cells.Item(16, 3).Value = 4047
$ws.Cells.Item(16, 4).Value = 5745197
$ws.Cells.Item(20, 3).Value = 7076
$ws.Cells.Item(20, 4).Value = 9881318
$ws.Cells.Item(22, 3).Value = 78986
$ws.Cells.Item(22, 4).Value = 98399880
$ws.Cells.Item(23, 3).Value = 56
$ws.Cells.Item(23, 4).Value = 74209
$ws.Cells.Item(28, 3).Value = 32926
$ws.Cells.Item(28, 4).Value = 48189214
$ws.Cells.Item(30, 3).Value = 11712
$ws.Cells.Item(30, 4).Value = 16848427
$ws.Cells.Item(35, 3).Value = 1945
$ws.Cells.Item(35, 4).Value = 2745194
$ws.Cells.Item(36, 3).Value = 98954
$ws.Cells.Item(36, 4).Value = 124421860
$ws.Cells.Item(42, 3).Value = 908
$ws.Cells.Item(42, 4).Value = 1336685
$ws.Cells.Item(44, 3).Value = 44965
$ws.Cells.Item(44, 4).Value = 65892848
$ws.Cells.Item(46, 3).Value = 9349
$ws.Cells.Item(46, 4).Value = 13408625
$ws.Cells.Item(48, 3).Value = 1426
$ws.Cells.Item(48, 4).Value = 1981311
$ws.Cells.Item(51, 3).Value = 2474
$ws.Cells.Item(51, 4).Value = 3459602
$ws.Cells.Item(52, 3).Value = 70406
$ws.Cells.Item(52, 4).Value = 88281022
$ws.Cells.Item(53, 3).Value = 44
$ws.Cells.Item(53, 4).Value = 50883
$ws.Cells.Item(59, 3).Value = 28627
$ws.Cells.Item(59, 4).Value = 41983766
$ws.Cells.Item(62, 3).Value = 11384
$ws.Cells.Item(62, 4).Value = 16458647
$ws.Cells.Item(64, 3).Value = 1373
$ws.Cells.Item(64, 4).Value = 1919497
$ws.Cells.Item(68, 3).Value = 1574
$ws.Cells.Item(68, 4).Value = 2206581
$ws.Cells.Item(70, 3).Value = 20885
$ws.Cells.Item(70, 4).Value = 27348823
$ws.Cells.Item(74, 3).Value = 7722
$ws.Cells.Item(74, 4).Value = 11309601
$ws.Cells.Item(76, 3).Value = 5214
$ws.Cells.Item(76, 4).Value = 7571499
$ws.Cells.Item(77, 3).Value = 497
$ws.Cells.Item(77, 4).Value = 704239
$ws.Cells.Item(78, 3).Value = 290
$ws.Cells.Item(78, 4).Value = 408083
$ws.Cells.Item(79, 3).Value = 143721
$ws.Cells.Item(79, 4).Value = 179093809
$ws.Cells.Item(80, 3).Value = 71
$ws.Cells.Item(80, 4).Value = 84766
$ws.Cells.Item(81, 3).Value = 89
$ws.Cells.Item(81, 4).Value = 124884
$ws.Cells.Item(83, 3).Value = 437
$ws.Cells.Item(83, 4).Value = 638324
$ws.Cells.Item(85, 3).Value = 64547
$ws.Cells.Item(85, 4).Value = 94593284
$ws.Cells.Item(88, 3).Value = 30321
$ws.Cells.Item(88, 4).Value = 43863090
$ws.Cells.Item(90, 3).Value = 2769
$ws.Cells.Item(90, 4).Value = 3987143
$ws.Cells.Item(91, 3).Value = 2988
$ws.Cells.Item(91, 4).Value = 4225334
$ws.Cells.Item(92, 3).Value = 34335
$ws.Cells.Item(92, 4).Value = 46550683
$ws.Cells.Item(96, 3).Value = 8320
$ws.Cells.Item(96, 4).Value = 12230763
$ws.Cells.Item(98, 3).Value = 7676
$ws.Cells.Item(98, 4).Value = 11138362
$ws.Cells.Item(100, 3).Value = 549
$ws.Cells.Item(100, 4).Value = 779406
$ws.Cells.Item(101, 3).Value = 519
$ws.Cells.Item(101, 4).Value = 749050
$ws.Cells.Item(102, 3).Value = 11249
$ws.Cells.Item(102, 4).Value = 17824918
$ws.Cells.Item(104, 3).Value = 2754
$ws.Cells.Item(104, 4).Value = 4684854
$ws.Cells.Item(106, 3).Value = 3777
$ws.Cells.Item(106, 4).Value = 6448776
$ws.Cells.Item(108, 3).Value = 166
$ws.Cells.Item(108, 4).Value = 280445
$ws.Cells.Item(109, 3).Value = 220
$ws.Cells.Item(109, 4).Value = 356030
$ws.Cells.Item(110, 3).Value = 144666
$ws.Cells.Item(110, 4).Value = 178919200
$ws.Cells.Item(116, 3).Value = 53576
$ws.Cells.Item(116, 4).Value = 78520350
$ws.Cells.Item(118, 3).Value = 27864
$ws.Cells.Item(118, 4).Value = 40372562
$ws.Cells.Item(119, 3).Value = 1325
$ws.Cells.Item(119, 4).Value = 1812194
$ws.Cells.Item(122, 3).Value = 2377
$ws.Cells.Item(122, 4).Value = 3341770
$ws.Cells.Item(124, 3).Value = 535398
$ws.Cells.Item(124, 4).Value = 707450560
$ws.Cells.Item(125, 3).Value = 93
$ws.Cells.Item(125, 4).Value = 123789
$ws.Cells.Item(129, 3).Value = 1403
$ws.Cells.Item(129, 4).Value = 2079682
$ws.Cells.Item(131, 3).Value = 213701
$ws.Cells.Item(131, 4).Value = 314130371
$ws.Cells.Item(132, 3).Value = 418
$ws.Cells.Item(132, 4).Value = 623710
$ws.Cells.Item(134, 3).Value = 191244
$ws.Cells.Item(134, 4).Value = 278116760
$ws.Cells.Item(137, 3).Value = 2886
$ws.Cells.Item(137, 4).Value = 4052272
$ws.Cells.Item(140, 3).Value = 6751
$ws.Cells.Item(140, 4).Value = 9528706
$ws.Cells.Item(143, 3).Value = 45963
$ws.Cells.Item(143, 4).Value = 61349800
$ws.Cells.Item(149, 3).Value = 14375
$ws.Cells.Item(149, 4).Value = 21070307
$ws.Cells.Item(150, 3).Value = 3858
$ws.Cells.Item(150, 4).Value = 5563092
$ws.Cells.Item(155, 3).Value = 413
$ws.Cells.Item(155, 4).Value = 582313
$ws.Cells.Item(156, 3).Value = 18111
$ws.Cells.Item(156, 4).Value = 23945906
$ws.Cells.Item(160, 3).Value = 7382
$ws.Cells.Item(160, 4).Value = 10742048
$ws.Cells.Item(162, 3).Value = 5157
$ws.Cells.Item(162, 4).Value = 7423474
$ws.Cells.Item(165, 3).Value = 278
$ws.Cells.Item(165, 4).Value = 397164
$ws.Cells.Item(167, 3).Value = 20587
$ws.Cells.Item(167, 4).Value = 35592563
$ws.Cells.Item(168, 3).Value = 2213
$ws.Cells.Item(168, 4).Value = 3830701
$ws.Cells.Item(169, 3).Value = 292
$ws.Cells.Item(169, 4).Value = 495089
$ws.Cells.Item(171, 3).Value = 71
$ws.Cells.Item(171, 4).Value = 131690
$ws.Cells.Item(172, 3).Value = 117
$ws.Cells.Item(172, 4).Value = 211449
$ws.Cells.Item(173, 3).Value = 89496
$ws.Cells.Item(173, 4).Value = 111791414
$ws.Cells.Item(180, 3).Value = 34407
$ws.Cells.Item(180, 4).Value = 50448058
$ws.Cells.Item(182, 3).Value = 13321
$ws.Cells.Item(182, 4).Value = 19247992
$ws.Cells.Item(184, 3).Value = 1264
$ws.Cells.Item(184, 4).Value = 1768933
$ws.Cells.Item(186, 3).Value = 1740
$ws.Cells.Item(186, 4).Value = 2441487
$ws.Cells.Item(188, 3).Value = 243132
$ws.Cells.Item(188, 4).Value = 302059535
$ws.Cells.Item(190, 3).Value = 173
$ws.Cells.Item(190, 4).Value = 249736
$ws.Cells.Item(194, 3).Value = 890
$ws.Cells.Item(194, 4).Value = 1309345
$ws.Cells.Item(196, 3).Value = 87873
$ws.Cells.Item(196, 4).Value = 128793334
$ws.Cells.Item(199, 3).Value = 33722
$ws.Cells.Item(199, 4).Value = 48549463
$ws.Cells.Item(202, 3).Value = 5182
$ws.Cells.Item(202, 4).Value = 7377358
$ws.Cells.Item(205, 3).Value = 5157
$ws.Cells.Item(205, 4).Value = 7140867
$ws.Cells.Item(208, 3).Value = 269503
$ws.Cells.Item(208, 4).Value = 333455676
$ws.Cells.Item(217, 3).Value = 96725
$ws.Cells.Item(217, 4).Value = 141493030
$ws.Cells.Item(218, 3).Value = 95
$ws.Cells.Item(218, 4).Value = 141699
$ws.Cells.Item(220, 3).Value = 52747
$ws.Cells.Item(220, 4).Value = 76233209
$ws.Cells.Item(223, 3).Value = 4732
$ws.Cells.Item(223, 4).Value = 6641035
$ws.Cells.Item(226, 3).Value = 6112
$ws.Cells.Item(226, 4).Value = 8466100
$ws.Cells.Item(228, 3).Value = 7
$ws.Cells.Item(228, 4).Value = 10500
$ws.Cells.Item(229, 3).Value = 108657
$ws.Cells.Item(229, 4).Value = 135816647
$ws.Cells.Item(232, 3).Value = 11
$ws.Cells.Item(232, 4).Value = 14147
$ws.Cells.Item(234, 3).Value = 570
$ws.Cells.Item(234, 4).Value = 832939
$ws.Cells.Item(236, 3).Value = 50251
$ws.Cells.Item(236, 4).Value = 73610986
$ws.Cells.Item(238, 3).Value = 12765
$ws.Cells.Item(238, 4).Value = 18359767
$ws.Cells.Item(240, 3).Value = 1909
$ws.Cells.Item(240, 4).Value = 2737382
$ws.Cells.Item(242, 3).Value = 2645
$ws.Cells.Item(242, 4).Value = 3704384
$ws.Cells.Item(243, 3).Value = 263879
$ws.Cells.Item(243, 4).Value = 333156164
$ws.Cells.Item(251, 3).Value = 97521
$ws.Cells.Item(251, 4).Value = 142890737
$ws.Cells.Item(254, 3).Value = 66868
$ws.Cells.Item(254, 4).Value = 96920073
$ws.Cells.Item(256, 3).Value = 2449
$ws.Cells.Item(256, 4).Value = 3454224
$ws.Cells.Item(259, 3).Value = 4892
$ws.Cells.Item(259, 4).Value = 6873294
